{"js": "// Replace the date line and each \"a\u00f7b=\" problem in the practice table with\n// the next day's values, in document order. The document body paragraphs\n// (including the one-run paragraphs that live inside each table cell) are\n// visited in order and each non-empty paragraph's text is swapped for the\n// corresponding new value below.\nconst replacements = [\n  \"2024-09-28 Saturday\",\n  \"20\u00f72=\",\n  \"44\u00f73=\",\n  \"10\u00f75=\",\n  \"20\u00f72=\",\n  \"71\u00f76=\",\n  \"26\u00f77=\",\n  \"77\u00f77=\",\n  \"82\u00f73=\",\n  \"10\u00f76=\",\n  \"79\u00f77=\",\n  \"71\u00f79=\",\n  \"56\u00f72=\",\n  \"91\u00f79=\",\n  \"53\u00f76=\",\n  \"34\u00f75=\",\n  \"64\u00f72=\",\n  \"11\u00f75=\",\n  \"57\u00f74=\",\n  \"60\u00f79=\",\n  \"86\u00f78=\",\n  \"22\u00f72=\",\n  \"38\u00f76=\",\n  \"61\u00f77=\",\n  \"88\u00f79=\",\n  \"13\u00f73=\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  // Skip the blank paragraphs that sit in the empty spacer table rows.\n  if (text === \"\") {\n    continue;\n  }\n  if (idx >= replacements.length) {\n    break;\n  }\n  para.insertText(replacements[idx], Word.InsertLocation.replace);\n  idx++;\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each \"a\u00f7b=\" problem in the practice table with\n# the next day's values, in document order. Word represents every table\n# cell (even the blank spacer rows) as one or more paragraphs, so we walk\n# $d.Paragraphs in order and only touch the ones that actually contain text\n# (trimming the trailing paragraph-mark / end-of-cell control characters\n# that Range.Text always reports).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    \"2024-09-28 Saturday\",\n    \"20\u00f72=\",\n    \"44\u00f73=\",\n    \"10\u00f75=\",\n    \"20\u00f72=\",\n    \"71\u00f76=\",\n    \"26\u00f77=\",\n    \"77\u00f77=\",\n    \"82\u00f73=\",\n    \"10\u00f76=\",\n    \"79\u00f77=\",\n    \"71\u00f79=\",\n    \"56\u00f72=\",\n    \"91\u00f79=\",\n    \"53\u00f76=\",\n    \"34\u00f75=\",\n    \"64\u00f72=\",\n    \"11\u00f75=\",\n    \"57\u00f74=\",\n    \"60\u00f79=\",\n    \"86\u00f78=\",\n    \"22\u00f72=\",\n    \"38\u00f76=\",\n    \"61\u00f77=\",\n    \"88\u00f79=\",\n    \"13\u00f73=\"\n)\n\n$idx = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $trimmed = $r.Text.TrimEnd([char]13, [char]7)\n    if ($trimmed.Length -gt 0) {\n        if ($idx -ge $replacements.Length) {\n            break\n        }\n        $r.Text = $replacements[$idx]\n        $idx = $idx + 1\n    }\n}\n"}
